$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Starting point: two sheets -> usersList, documentsList
# Target:         four sheets -> rolesList, usersList, documentsList, documentMetaDataList
# ---------------------------------------------------------------------------

# 1) Insert a new "rolesList" sheet in front of "usersList" (becomes the first tab)
$usersList = $wb.Worksheets.Item("usersList")
$rolesList = $wb.Worksheets.Add($usersList)
$rolesList.Name = "rolesList"

$rolesList.Range("A1").Value = "ROLE_NAME"
$rolesList.Range("A2").Value = "ADMIN"
$rolesList.Range("A3").Value = "USER"
$rolesList.Columns.Item(1).ColumnWidth = 11.85546875
$rolesList.Range("L28").Select()

# 2) Append a new "documentMetaDataList" sheet after "documentsList" (becomes the last tab)
#    (re-fetch the sheet reference - the sheet collection just changed above)
$documentsList = $wb.Worksheets.Item("documentsList")
$metaList = $wb.Worksheets.Add($null, $documentsList)
$metaList.Name = "documentMetaDataList"

$metaList.Range("A1").Value = "USER_NAME"
$metaList.Range("B1").Value = "DOCUMENT_NAME"
$metaList.Range("C1").Value = "META_DATA_NAME"
$metaList.Range("D1").Value = "META_DATA_VALUE"

$metaList.Range("A2").Value = "mariusz"
$metaList.Range("B2").Value = "testDocumentMariusz1.docx"
$metaList.Range("C2").Value = "metaDataName1"
$metaList.Range("D2").Value = "metaDataValue1"

$metaList.Range("A3").Value = "mariusz"
$metaList.Range("B3").Value = "testDocumentMariusz1.docx"
$metaList.Range("C3").Value = "metaDataName2"
$metaList.Range("D3").Value = "metaDataValue2"

$metaList.Range("A4").Value = "mariusz"
$metaList.Range("B4").Value = "testDocumentMariusz1.docx"
$metaList.Range("C4").Value = "metaDataName3"
$metaList.Range("D4").Value = "metaDataValue3"

$metaList.Range("A5").Value = "mariusz"
$metaList.Range("B5").Value = "testDocumentMariusz1.docx"
$metaList.Range("C5").Value = "metaDataName4"
$metaList.Range("D5").Value = "metaDataValue4"

$metaList.Columns.Item(1).ColumnWidth = 11.85546875
$metaList.Columns.Item(2).ColumnWidth = 26.7109375
$metaList.Columns.Item(3).ColumnWidth = 18.5703125
$metaList.Columns.Item(4).ColumnWidth = 18.85546875
$metaList.Range("I13").Select()

# 3) Update "documentsList" PATH_TO_DOCUMENT values + widen column B + change selection
#    (re-fetch again - the sheet collection changed again above)
$documentsList = $wb.Worksheets.Item("documentsList")
$documentsList.Range("B2").Value = "testData/documents/testDocumentMariusz1.docx"
$documentsList.Range("B3").Value = "testData/documents/testDocumentMariusz2.docx"
$documentsList.Range("B4").Value = "testData/documents/testDocumentUser123.docx"
$documentsList.Range("B5").Value = "testData/documents/testDocumentUser1234.pdf"
$documentsList.Columns.Item(2).ColumnWidth = 46.42578125

# 4) documentsList stays the active/selected tab, with an updated selected cell
$documentsList.Activate()
$documentsList.Range("E7").Select()
